# Commit: "Search Product check Table (cot ten)"
# The author was testing the product/category search feature and, while
# poking around the workbook, retyped a couple of lookup-table labels and
# left the selection in different cells on three of the sheets.

$wb = $excel.ActiveWorkbook

# --- Category sheet -------------------------------------------------
# The "ParentCategory" column (B) for every row now reads "Category 1"
# instead of the old placeholder brand-ish names.
$wsCategory = $wb.Worksheets.Item("Category")
$wsCategory.Range("B2").Value = "Category 1"
$wsCategory.Range("B3").Value = "Category 1"
$wsCategory.Range("B4").Value = "Category 1"

$wsCategory.Activate()
$wsCategory.Range("B2").Select()

# --- Brand sheet ------------------------------------------------------
# No data changed here, only where the user left the selection.
$wsBrand = $wb.Worksheets.Item("Brand")
$wsBrand.Activate()
$wsBrand.Range("C1").Select()

# --- Product sheet ------------------------------------------------------
# Category Name column (B) values corrected for the first two products.
$wsProduct = $wb.Worksheets.Item("Product")
$wsProduct.Range("B2").Value = "Mỹ Phẩm"
$wsProduct.Range("B3").Value = "Cake"

$wsProduct.Activate()
$wsProduct.Range("I2").Select()
